# Fruta / hortaliza, semanal
# Insert a new weekly record at row 125 (pushing the existing rows 125-140
# down to 126-141) on the Alcachofa - Macroferia Regional de Talca sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 125:140 down to 126:141, leaving a blank row 125 for the new record.
$ws.Rows("125:125").Insert()

# Populate the newly inserted row 125 with the new weekly data point.
$ws.Range("A125").Value = 5
$ws.Range("B125").Value = "Macroferia Regional de Talca"
$ws.Range("C125").Value = "Maule"
$ws.Range("D125").Value = 45166
$ws.Range("E125").Value = 7
$ws.Range("F125").Value = 100112013
$ws.Range("G125").Value = "Alcachofa"
$ws.Range("H125").Value = "Madrigal"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 13000
$ws.Range("L125").Value = 13000
$ws.Range("M125").Value = 13000
$ws.Range("N125").Value = "$/caja 40 unidades"
$ws.Range("O125").Value = "Provincia del Elquí"
$ws.Range("P125").Value = 325
$ws.Range("Q125").Value = 40
$ws.Range("R125").Value = "Hortaliza"
